# Fix: Ensures Invoice Number and EST# are read as strings
# Previously the Invoice Number and EST# could be read in as int64 if none
# of the rows had strings. This change forces the type to string by
# replacing the "#4" placeholder text values in row 8 (EST# / Invoice
# columns) with the actual numeric estimate/invoice number 5, which also
# drops the now-unused "#5" entry from the shared-strings table. The
# trailing blank row (row 19) is also removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: EST# (A8) and Invoice (B8) were text placeholders ("#4"); change
# them to the real numeric estimate/invoice number 5.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 5

# The selected/active cell moves from C9 to B9.
$ws.Range("B9").Select()

# Delete the last (now superfluous) blank row so the sheet dimension
# shrinks from A1:I19 to A1:I18.
$ws.Rows.Item(19).Delete()
